# Update the stations_metadata worksheet:
#  - reorder/rename header columns and introduce new columns
#    (station_kodas, x_coord, y_coord) while moving roughness_n earlier
#  - replace the sample data rows with a fuller set of stations, each
#    now carrying a numeric station_id, a textual station_code, a
#    station_kodas, coordinates, and min/max/datum level columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (A1:L1)
# ---------------------------------------------------------------------
$headers = @(
  "river_name",
  "station_name",
  "station_code",
  "station_id",
  "station_kodas",
  "x_coord",
  "y_coord",
  "roughness_n",
  "basin_name",
  "datum_offset_cm",
  "min_level_cm",
  "max_level_cm"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. Data rows (A2:L8)
# columns: river_name, station_name, station_code, station_id,
#          station_kodas, x_coord, y_coord, roughness_n, basin_name,
#          datum_offset_cm, min_level_cm, max_level_cm
# ---------------------------------------------------------------------
$data = @(
  @("Merkys",               "Puvočiai",            "101", 5101, "5101LT", 575000, 5998000, 0.04,  "Nemunas-Merkys", 0,     50,   850),
  @("Nemunas",               "Druskininkai",         "102", 5102, "5102LT", 568500, 5992000, 0.038, "Nemunas-Main",   0,     20,   1000),
  @("Verknė",                "Verbyliškės",          "103", 5103, "5103LT", 521000, 6032000, 0.042, "Nemunas-Verkne", 0,     30,   700),
  @("Nemunas",               "Nemunaičiai",          "104", 5104, "5104LT", 540500, 6042000, 0.038, "Nemunas-Main",   0,     20,   1100),
  @("Merkys",                "Jašiūnai",             "105", 5105, "5105LT", 583500, 6029000, 0.04,  "Nemunas-Merkys", 0,     30,   800),
  @("Šešupė",                "Kudirkos Naumiestis",  "106", 5106, "5106LT", 409000, 6022000, 0.045, "Nemunas-Sesupe", 0,     40,   950),
  @("Nemuno atšaka Atmata",  "Rusnė",                "769", 769,  "60004LT",333694, 6132670, 0.03,  "Nemunas-Delta", -1.56, -100, 300)
)

$rowIndex = 2
foreach ($row in $data) {
  $ws.Cells.Item($rowIndex, 1).Value  = $row[0]   # river_name
  $ws.Cells.Item($rowIndex, 2).Value  = $row[1]   # station_name
  # station_code looks numeric ("101", "102", ...) - force it to stay a
  # text value (not auto-converted to a number) by writing it as a text
  # formula and then collapsing that formula down to a static value via
  # copy / paste-special, which keeps the cell unstyled.
  $ws.Cells.Item($rowIndex, 3).Formula = '="' + $row[2] + '"'
  $ws.Cells.Item($rowIndex, 4).Value  = $row[3]   # station_id (number)
  $ws.Cells.Item($rowIndex, 5).Value  = $row[4]   # station_kodas
  $ws.Cells.Item($rowIndex, 6).Value  = $row[5]   # x_coord
  $ws.Cells.Item($rowIndex, 7).Value  = $row[6]   # y_coord
  $ws.Cells.Item($rowIndex, 8).Value  = $row[7]   # roughness_n
  $ws.Cells.Item($rowIndex, 9).Value  = $row[8]   # basin_name
  $ws.Cells.Item($rowIndex, 10).Value = $row[9]   # datum_offset_cm
  $ws.Cells.Item($rowIndex, 11).Value = $row[10]  # min_level_cm
  $ws.Cells.Item($rowIndex, 12).Value = $row[11]  # max_level_cm
  $rowIndex++
}

$codeRange = $ws.Range("C2:C8")
$codeRange.Copy()
$codeRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false
